$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 746.5
$ws.Range("I2").Value = 747
$ws.Range("J2").Value = 746
$ws.Range("K2").Value = 747
$ws.Range("L2").Value = 746
$ws.Range("M2").Value = -634
$ws.Range("N2").Value = -972
$ws.Range("H4").Value = 222.46153
$ws.Range("I4").Value = 88.55556
$ws.Range("J4").Value = 523.75
$ws.Range("K4").Value = 88.55556
$ws.Range("L4").Value = 523.75
$ws.Range("M4").Value = 25.44444
$ws.Range("N4").Value = -751.75
$ws.Range("H18").Value = 932.8333
$ws.Range("I18").Value = 932.8333
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 932.8333
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -648.8333
$ws.Range("N18").ClearContents()
$ws.Range("H33").Value = 430.53333
$ws.Range("I33").Value = 390
$ws.Range("K33").Value = 390
$ws.Range("M33").Value = -161
$ws.Range("H95").Value = 39966.332
$ws.Range("J95").Value = 39966.332
$ws.Range("L95").Value = 39966.332
$ws.Range("N95").Value = -45458.332
$ws.Range("H98").Value = 490.7
$ws.Range("I98").Value = 378.77777
$ws.Range("K98").Value = 378.77777
$ws.Range("M98").Value = 1119.22223
$ws.Range("H106").Value = 1249.25
$ws.Range("I106").Value = 999
$ws.Range("K106").Value = 999
$ws.Range("M106").Value = -368
$ws.Range("H113").Value = 8000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -1746
$ws.Range("H122").Value = 490.7
$ws.Range("I122").Value = 378.77777
$ws.Range("K122").Value = 1136.33331
$ws.Range("M122").Value = 1313.66669
$ws.Range("H131").Value = 4729.5
$ws.Range("I131").Value = 490
$ws.Range("K131").Value = 1470
$ws.Range("M131").Value = 3570
$ws.Range("H132").Value = 1150
$ws.Range("I132").Value = 1150
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3450
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -920
$ws.Range("N132").ClearContents()
$ws.Range("H138").Value = 2691.8
$ws.Range("I138").Value = 2395.818
$ws.Range("J138").Value = 3505.75
$ws.Range("K138").Value = 7187.454000000001
$ws.Range("L138").Value = 10517.25
$ws.Range("M138").Value = -2047.454000000001
$ws.Range("N138").Value = -20797.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1794.1428
$ws.Range("I32").Value = 1794.1428
$ws.Range("K32").Value = 1794.1428
$ws.Range("M32").Value = -1507.1428
$ws.Range("H63").Value = 4159.4
$ws.Range("I63").Value = 3266.6667
$ws.Range("J63").Value = 5498.5
$ws.Range("K63").Value = 3266.6667
$ws.Range("L63").Value = 5498.5
$ws.Range("M63").Value = -2580.6667
$ws.Range("N63").Value = -6870.5
$ws.Range("H66").Value = 4159.4
$ws.Range("I66").Value = 3266.6667
$ws.Range("J66").Value = 5498.5
$ws.Range("K66").Value = 16333.3335
$ws.Range("L66").Value = 27492.5
$ws.Range("M66").Value = -12901.3335
$ws.Range("N66").Value = -34356.5
$ws.Range("H88").Value = 2250
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 2250
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4399.75
$ws.Range("I99").Value = 4399.75
$ws.Range("K99").Value = 4399.75
$ws.Range("M99").Value = -2901.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 661
$ws.Range("I22").Value = 562.8182
$ws.Range("K22").Value = 562.8182
$ws.Range("M22").Value = -212.8182
$ws.Range("H31").Value = 3357.2
$ws.Range("I31").Value = 3357.2
$ws.Range("K31").Value = 3357.2
$ws.Range("M31").Value = -3062.2
$ws.Range("H34").Value = 3357.2
$ws.Range("I34").Value = 3357.2
$ws.Range("K34").Value = 3357.2
$ws.Range("M34").Value = -3155.2
$ws.Range("H93").Value = 32333.334
$ws.Range("I93").Value = 17500
$ws.Range("J93").Value = 62000
$ws.Range("K93").Value = 17500
$ws.Range("L93").Value = 62000
$ws.Range("M93").Value = -15628
$ws.Range("N93").Value = -65744
$ws.Range("H107").Value = 698.13336
$ws.Range("I107").Value = 534.4545000000001
$ws.Range("J107").Value = 1148.25
$ws.Range("K107").Value = 534.4545000000001
$ws.Range("L107").Value = 1148.25
$ws.Range("M107").Value = 1385.5455
$ws.Range("N107").Value = -4988.25
$ws.Range("H132").Value = 5963.9375
$ws.Range("I132").Value = 5963.9375
$ws.Range("K132").Value = 17891.8125
$ws.Range("M132").Value = -15361.8125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 10000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 30000
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -30168
$ws.Range("H52").Value = 2000
$ws.Range("J52").Value = 2000
$ws.Range("L52").Value = 6000
$ws.Range("N52").Value = -6532
$ws.Range("H60").Value = 598.6667
$ws.Range("I60").Value = 598.6667
$ws.Range("K60").Value = 1796.0001
$ws.Range("M60").Value = -1545.0001
$ws.Range("H81").Value = 19498.334
$ws.Range("J81").Value = 28747.5
$ws.Range("L81").Value = 86242.5
$ws.Range("N81").Value = -88488.5
$ws.Range("H84").Value = 19498.334
$ws.Range("J84").Value = 28747.5
$ws.Range("L84").Value = 258727.5
$ws.Range("N84").Value = -269959.5
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("N115").ClearContents()
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H139").Value = 1766.5
$ws.Range("I139").Value = 1745.2727
$ws.Range("J139").Value = 2000
$ws.Range("K139").Value = 5235.8181
$ws.Range("L139").Value = 6000
$ws.Range("M139").Value = -95.81810000000041
$ws.Range("N139").Value = -16280

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1000
$ws.Range("I80").Value = 1000
$ws.Range("K80").Value = 1000
$ws.Range("M80").Value = -2
$ws.Range("H83").Value = 1000
$ws.Range("I83").Value = 1000
$ws.Range("K83").Value = 5000
$ws.Range("M83").Value = -8
$ws.Range("H124").Value = 100000
$ws.Range("J124").Value = 100000
$ws.Range("L124").Value = 100000
$ws.Range("N124").Value = -109820

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 13382.429
$ws.Range("I40").Value = 11950.363
$ws.Range("K40").Value = 11950.363
$ws.Range("M40").Value = -11814.363
$ws.Range("H46").Value = 482.33334
$ws.Range("I46").Value = 481.56522
$ws.Range("K46").Value = 481.56522
$ws.Range("M46").Value = -293.56522
$ws.Range("H53").Value = 9000
$ws.Range("I53").Value = 10000
$ws.Range("J53").Value = 7000
$ws.Range("K53").Value = 10000
$ws.Range("L53").Value = 7000
$ws.Range("M53").Value = -9482
$ws.Range("N53").Value = -8036
$ws.Range("H68").Value = 4800
$ws.Range("J68").Value = 4800
$ws.Range("L68").Value = 4800
$ws.Range("N68").Value = -6298
$ws.Range("H71").Value = 4800
$ws.Range("J71").Value = 4800
$ws.Range("L71").Value = 24000
$ws.Range("N71").Value = -31488
$ws.Range("H82").Value = 4162.5
$ws.Range("I82").Value = 4000
$ws.Range("J82").Value = 4216.6665
$ws.Range("K82").Value = 4000
$ws.Range("L82").Value = 4216.6665
$ws.Range("M82").Value = -3639
$ws.Range("N82").Value = -4938.6665
$ws.Range("H85").Value = 4162.5
$ws.Range("I85").Value = 4000
$ws.Range("J85").Value = 4216.6665
$ws.Range("K85").Value = 4000
$ws.Range("L85").Value = 4216.6665
$ws.Range("M85").Value = -2752
$ws.Range("N85").Value = -6712.6665
$ws.Range("H88").Value = 10000
$ws.Range("I88").Value = 10000
$ws.Range("K88").Value = 10000
$ws.Range("M88").Value = -9572
$ws.Range("H91").Value = 10000
$ws.Range("I91").Value = 10000
$ws.Range("K91").Value = 10000
$ws.Range("M91").Value = -8518

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2511999.8
$ws.Range("I5").Value = 10000000
$ws.Range("J5").Value = 639999.75
$ws.Range("K5").Value = 10000000
$ws.Range("L5").Value = 639999.75
$ws.Range("M5").Value = -9999888
$ws.Range("N5").Value = -640223.75
$ws.Range("H62").Value = 2008.3334
$ws.Range("I62").Value = 1250
$ws.Range("J62").Value = 2387.5
$ws.Range("K62").Value = 1250
$ws.Range("L62").Value = 2387.5
$ws.Range("M62").Value = -626
$ws.Range("N62").Value = -3635.5
$ws.Range("H65").Value = 2008.3334
$ws.Range("I65").Value = 1250
$ws.Range("J65").Value = 2387.5
$ws.Range("K65").Value = 6250
$ws.Range("L65").Value = 11937.5
$ws.Range("M65").Value = -3130
$ws.Range("N65").Value = -18177.5
$ws.Range("H96").Value = 1367.4
$ws.Range("I96").Value = 614.8
$ws.Range("J96").Value = 2120
$ws.Range("K96").Value = 614.8
$ws.Range("L96").Value = 2120
$ws.Range("M96").Value = 758.2
$ws.Range("N96").Value = -4866
